# "correcting coordinate system image again"
#
# Repositions several of the red/green/blue coordinate-axis arrow
# connectors and their axis-label textboxes on the single slide of the
# coordinate-systems deck, and moves one connector ("Gerade Verbindung
# mit Pfeil 10", id 11) to the end of the shape stack (so it renders on
# top of everything else) while un-flipping it horizontally.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $cand = $shapes.Item($i)
        if ($cand.Id -eq $id) {
            return $cand
        }
    }
    return $null
}

# --- Shape id 5: "Gerade Verbindung mit Pfeil 4" (red y-axis arrow) ---
$shape5 = Get-ShapeById $s.Shapes 5
$shape5.Left   = 435.7649841308594
$shape5.Top    = 286.5570373535156
$shape5.Width  = 4.2853546142578125
$shape5.Height = 61.109134674072266

# --- Shape id 10: "Textfeld 9" (the "x" label near the first axes) ---
$shape10 = Get-ShapeById $s.Shapes 10
$shape10.Left = 438.03033447265625
$shape10.Top  = 325.4737243652344

# --- Shape id 11: "Gerade Verbindung mit Pfeil 10" (red arrow, second axes) ---
# Moves to a new spot, loses its horizontal flip, and gets sent to the very
# end of the shape tree (front-most / top of z-order).
$shape11 = Get-ShapeById $s.Shapes 11
$shape11.HorizontalFlip = $false
$shape11.Left   = 587.2384643554688
$shape11.Top    = 242.5303955078125
$shape11.Width  = 6.018425464630127
$shape11.Height = 48.32448959350586
$shape11.ZOrder(0)   # msoBringToFront -> becomes the last shape in spTree

# --- Shape id 12: "Gerade Verbindung mit Pfeil 11" (blue arrow, second axes) ---
$shape12 = Get-ShapeById $s.Shapes 12
$shape12.Left   = 580.299072265625
$shape12.Top    = 241.6519775390625
$shape12.Width  = 7.245197296142578
$shape12.Height = 56.681419372558594

# --- Shape id 16: "Textfeld 15" (the "z" label) ---
$shape16 = Get-ShapeById $s.Shapes 16
$shape16.Left = 568.3798828125
$shape16.Top  = 294.94561767578125

# --- Shape id 17: "Textfeld 16" (the "x" label, second axes) ---
$shape17 = Get-ShapeById $s.Shapes 17
$shape17.Left = 586.6693725585938
$shape17.Top  = 283.75347900390625
